$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.765.90"
$ws.Range("E2").Value = "  +1.26%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.337.78"
$ws.Range("E3").Value = "  +2.00%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.79"
$ws.Range("E5").Value = "  +0.57%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "175.82"
$ws.Range("E6").Value = "  +2.22%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("E8").Value = "  +1.95%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.331.09"
$ws.Range("E9").Value = "  +1.80%  "
$ws.Range("E10").Value = "  +6.39%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.580"
$ws.Range("E11").Value = "  +1.71%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "46.91"
$ws.Range("E12").Value = "  +5.09%  "
$ws.Range("E13").Value = "  +1.51%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "688.75"
$ws.Range("E14").Value = "  +1.04%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.883.71"
$ws.Range("E15").Value = "  +2.26%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "8.42"
$ws.Range("E16").Value = "  +2.51%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.804.90"
$ws.Range("E17").Value = "  +1.17%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.118"
$ws.Range("E18").Value = "  +0.23%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.351.69"
$ws.Range("E19").Value = "  +2.53%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.58"
$ws.Range("E20").Value = "  +2.52%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.05"
$ws.Range("E21").Value = "  +4.28%  "
$ws.Range("E22").Value = "  +1.62%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.43"
$ws.Range("E23").Value = "  +3.69%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "16.93"
$ws.Range("E24").Value = "  +0.85%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "101.64"
$ws.Range("E25").Value = "  +3.26%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.90"
$ws.Range("E26").Value = "  +2.23%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.68"
$ws.Range("E27").Value = "  +2.05%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.49"
$ws.Range("E28").Value = "  +5.95%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.01"
$ws.Range("E29").Value = "  +0.10%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.55"
$ws.Range("E30").Value = "  +3.72%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.05"
$ws.Range("E31").Value = "  +7.52%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "568.73"
$ws.Range("E32").Value = "  -1.40%  "
$ws.Range("E33").Value = "  +1.94%  "
$ws.Range("E34").Value = "  +3.24%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "57.36"
$ws.Range("E35").Value = "  +3.59%  "
$ws.Range("B36").Value = "Maker"
$ws.Range("C36").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.717.30"
$ws.Range("E36").Value = "  -2.11%  "
$ws.Range("B37").Value = "Dai"
$ws.Range("C37").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("E37").Value = "  +0.00%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.29"
$ws.Range("E38").Value = "  +1.35%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "34.96"
$ws.Range("E39").Value = "  +12.00%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.134"
$ws.Range("E40").Value = "  +5.53%  "
$ws.Range("E41").Value = "  +7.16%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.62"
$ws.Range("E42").Value = "  +2.92%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0₃0674"
$ws.Range("E43").Value = "  +3.62%  "
$ws.Range("B44").Value = "ApeXProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.32"
$ws.Range("E44").Value = "  -2.10%  "
$ws.Range("B45").Value = "TheGraph"
$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.335"
$ws.Range("E45").Value = "  +3.82%  "
$ws.Range("E46").Value = "  +2.54%  "
$ws.Range("E47").Value = "  +6.05%  "
$ws.Range("E48").Value = "  +1.96%  "
$ws.Range("E49").Value = "  +0.05%  "
$ws.Range("E50").Value = "  +0.98%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "132.06"
$ws.Range("E51").Value = "  +3.47%  "
